$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Insert the new "Vehicle Model" row as row 54 (pushes the old row 54
#    "Vehicle Primary Color" down to row 55, and everything below shifts by
#    one). Copy row 53's formatting down first so the inserted row picks up
#    the same cell styles (s="3" for A-D, s="11" for E) as its neighbours.
# ---------------------------------------------------------------------------
$ws.Rows.Item(53).Copy()
$ws.Rows.Item(54).Insert()

$ws.Range("B54").Value = "Model of the vehicle"
$ws.Range("C54").Value = "Model"
$ws.Range("D54").Value = "Corolla"
$ws.Range("E54").Value = "/ir-doc:IncidentReport/lexspd:doPublish/lexs:PublishMessageContainer/lexs:PublishMessage/lexs:DataItemPackage/lexs:Digest/lexsdigest:EntityVehicle/nc:Vehicle[@s:id=/ir-doc:IncidentReport/lexspd:doPublish/lexs:PublishMessageContainer/lexs:PublishMessage/lexs:DataItemPackage/lexs:Digest/lexsdigest:Associations/nc:PersonConveyanceAssociation[nc:PersonReference/@s:ref=/ir-doc:IncidentReport/lexspd:doPublish/lexs:PublishMessageContainer/lexs:PublishMessage/lexs:DataItemPackage/lexs:Digest/lexsdigest:EntityPerson[j:CitationSubject]/lexsdigest:Person/@s:id]/nc:ConveyanceReference/@s:ref]/nc:ItemModelName"
$ws.Range("A54").Value = "Vehicle Model"

# ---------------------------------------------------------------------------
# 2) Insert the new "Officer Notes" row. After the insertion above, it needs
#    to land at row 104 (immediately above the old row 103 "Traffic Stop
#    Reason Code", which is now row 104 before this second insert).
#    Copy row 53's formatting again (A/D -> s="3", E -> s="11") and then
#    clear the B/C cells entirely so they don't materialize in the row
#    (the target row only has A, D and E populated).
# ---------------------------------------------------------------------------
$ws.Rows.Item(53).Copy()
$ws.Rows.Item(104).Insert()

$ws.Range("B104").Clear()
$ws.Range("C104").Clear()

$ws.Range("A104").Value = "Officer Notes"
$ws.Range("D104").Value = "Officer Notes"
$ws.Range("E104").Value = "/ir-doc:IncidentReport/lexspd:doPublish/lexs:PublishMessageContainer/lexs:PublishMessage/lexs:DataItemPackage/lexs:StructuredPayload/inc-ext:IncidentReport/inc-ext:DrivingIncident[lexslib:SameAsDigestReference/@lexslib:ref=/ir-doc:IncidentReport/lexspd:doPublish/lexs:PublishMessageContainer/lexs:PublishMessage/lexs:DataItemPackage/lexs:Digest/lexsdigest:EntityDocument/nc:Document/@s:id]/inc-ext:EnforcementOfficialNotesText"

# ---------------------------------------------------------------------------
# 3) Update the view/selection so the active cell matches the author's final
#    position.
# ---------------------------------------------------------------------------
$ws.Range("B103").Select() | Out-Null
